$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 345.25
$ws.Range("I2").Value = 345.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 345.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -232.25

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H19").Value = 3993.4
$ws.Range("I19").Value = 3624
$ws.Range("J19").Value = 4239.6665
$ws.Range("K19").Value = 3624
$ws.Range("L19").Value = 4239.6665
$ws.Range("M19").Value = -3449
$ws.Range("N19").Value = -4589.6665

$ws.Range("H47").Value = 20000
$ws.Range("I47").Value = 20000
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 20000
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -19028

$ws.Range("H70").Value = 4559.6665
$ws.Range("I70").Value = 2467.3333
$ws.Range("J70").Value = 6652
$ws.Range("K70").Value = 7401.999899999999
$ws.Range("L70").Value = 19956
$ws.Range("M70").Value = -7131.999899999999
$ws.Range("N70").Value = -20496

$ws.Range("H73").Value = 4559.6665
$ws.Range("I73").Value = 2467.3333
$ws.Range("J73").Value = 6652
$ws.Range("K73").Value = 7401.999899999999
$ws.Range("L73").Value = 19956
$ws.Range("M73").Value = -6465.999899999999
$ws.Range("N73").Value = -21828

$ws.Range("H86").Value = 11553
$ws.Range("I86").Value = 3121
$ws.Range("J86").Value = 19985
$ws.Range("K86").Value = 3121
$ws.Range("L86").Value = 19985
$ws.Range("M86").Value = -1998
$ws.Range("N86").Value = -22231

$ws.Range("H89").Value = 11553
$ws.Range("I89").Value = 3121
$ws.Range("J89").Value = 19985
$ws.Range("K89").Value = 15605
$ws.Range("L89").Value = 99925
$ws.Range("M89").Value = -9989
$ws.Range("N89").Value = -111157

$ws.Range("H137").Value = 1224.75
$ws.Range("I137").Value = 1224.75
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3674.25
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1124.25

$ws.Range("H141").Value = 4611.2666
$ws.Range("I141").Value = 4024.6365
$ws.Range("J141").Value = 6224.5
$ws.Range("K141").Value = 12073.9095
$ws.Range("L141").Value = 18673.5
$ws.Range("M141").Value = -6893.9095
$ws.Range("N141").Value = -29033.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1895.7142
$ws.Range("I63").Value = 1873.625
$ws.Range("J63").Value = 1909.3077
$ws.Range("K63").Value = 1873.625
$ws.Range("L63").Value = 1909.3077
$ws.Range("M63").Value = -1187.625
$ws.Range("N63").Value = -3281.3077

$ws.Range("H66").Value = 1895.7142
$ws.Range("I66").Value = 1873.625
$ws.Range("J66").Value = 1909.3077
$ws.Range("K66").Value = 9368.125
$ws.Range("L66").Value = 9546.538500000001
$ws.Range("M66").Value = -5936.125
$ws.Range("N66").Value = -16410.5385

$ws.Range("H122").Value = 860.2
$ws.Range("I122").Value = 651.5
$ws.Range("J122").Value = 1695
$ws.Range("K122").Value = 1954.5
$ws.Range("L122").Value = 5085
$ws.Range("M122").Value = 495.5
$ws.Range("N122").Value = -9985

$ws.Range("H132").Value = 15156825
$ws.Range("I132").Value = 3807.64
$ws.Range("J132").Value = 62510004
$ws.Range("K132").Value = 11422.92
$ws.Range("L132").Value = 187530012
$ws.Range("M132").Value = -8892.92
$ws.Range("N132").Value = -187535072

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2315.6
$ws.Range("I20").Value = 1794.5
$ws.Range("J20").Value = 4400
$ws.Range("K20").Value = 1794.5
$ws.Range("L20").Value = 4400
$ws.Range("M20").Value = -1547.5
$ws.Range("N20").Value = -4894

$ws.Range("H99").Value = 1437.909
$ws.Range("I99").Value = 1140.5
$ws.Range("J99").Value = 1794.8
$ws.Range("K99").Value = 1140.5
$ws.Range("L99").Value = 1794.8
$ws.Range("M99").Value = 357.5
$ws.Range("N99").Value = -4790.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 25004996
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 25004996
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 25004996
$ws.Range("N4").Value = -25005220

$ws.Range("H22").Value = 6742.091
$ws.Range("I22").Value = 50000
$ws.Range("J22").Value = 2416.3
$ws.Range("K22").Value = 50000
$ws.Range("L22").Value = 2416.3
$ws.Range("M22").Value = -49650
$ws.Range("N22").Value = -3116.3

$ws.Range("H86").Value = 3498.5
$ws.Range("I86").Value = 3498
$ws.Range("J86").Value = 3499
$ws.Range("K86").Value = 3498
$ws.Range("L86").Value = 3499
$ws.Range("M86").Value = -2375
$ws.Range("N86").Value = -5745

$ws.Range("H89").Value = 3498.5
$ws.Range("I89").Value = 3498
$ws.Range("J89").Value = 3499
$ws.Range("K89").Value = 17490
$ws.Range("L89").Value = 17495
$ws.Range("M89").Value = -11874
$ws.Range("N89").Value = -28727

$ws.Range("H99").Value = 2384.1667
$ws.Range("I99").Value = 2205.5
$ws.Range("J99").Value = 2473.5
$ws.Range("K99").Value = 2205.5
$ws.Range("L99").Value = 2473.5
$ws.Range("M99").Value = -707.5
$ws.Range("N99").Value = -5469.5

$ws.Range("H107").Value = 4935.75
$ws.Range("I107").Value = 3393.8
$ws.Range("J107").Value = 7505.6665
$ws.Range("K107").Value = 3393.8
$ws.Range("L107").Value = 7505.6665
$ws.Range("M107").Value = -1473.8
$ws.Range("N107").Value = -11345.6665

$ws.Range("H126").Value = 2384.1667
$ws.Range("I126").Value = 2205.5
$ws.Range("J126").Value = 2473.5
$ws.Range("K126").Value = 6616.5
$ws.Range("L126").Value = 7420.5
$ws.Range("M126").Value = -4146.5
$ws.Range("N126").Value = -12360.5

$ws.Range("H134").Value = 2253.4375
$ws.Range("I134").Value = 2237
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 6711
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -4176
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4609.1113
$ws.Range("I80").Value = 4349.5
$ws.Range("J80").Value = 4683.2856
$ws.Range("K80").Value = 13048.5
$ws.Range("L80").Value = 14049.8568
$ws.Range("M80").Value = -12112.5
$ws.Range("N80").Value = -15921.8568

$ws.Range("H83").Value = 4609.1113
$ws.Range("I83").Value = 4349.5
$ws.Range("J83").Value = 4683.2856
$ws.Range("K83").Value = 39145.5
$ws.Range("L83").Value = 42149.5704
$ws.Range("M83").Value = -34465.5
$ws.Range("N83").Value = -51509.5704

$ws.Range("H116").Value = 106451.7
$ws.Range("I116").Value = 117279.664
$ws.Range("J116").Value = 9000
$ws.Range("K116").Value = 351838.992
$ws.Range("L116").Value = 27000
$ws.Range("M116").Value = -348396.992
$ws.Range("N116").Value = -33884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3884.889
$ws.Range("I80").Value = 2369.5
$ws.Range("J80").Value = 5097.2
$ws.Range("K80").Value = 2369.5
$ws.Range("L80").Value = 5097.2
$ws.Range("M80").Value = -1371.5
$ws.Range("N80").Value = -7093.2

$ws.Range("H83").Value = 3884.889
$ws.Range("I83").Value = 2369.5
$ws.Range("J83").Value = 5097.2
$ws.Range("K83").Value = 11847.5
$ws.Range("L83").Value = 25486
$ws.Range("M83").Value = -6855.5
$ws.Range("N83").Value = -35470

$ws.Range("H97").Value = 1164.4286
$ws.Range("I97").Value = 1164.4286
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1164.4286
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -668.4286
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 2470.3215
$ws.Range("I102").Value = 1886.5294
$ws.Range("J102").Value = 3372.5454
$ws.Range("K102").Value = 1886.5294
$ws.Range("L102").Value = 3372.5454
$ws.Range("M102").Value = -264.5293999999999
$ws.Range("N102").Value = -6616.5454

$ws.Range("H126").Value = 2579.1428
$ws.Range("I126").Value = 2003.3334
$ws.Range("J126").Value = 2736.182
$ws.Range("K126").Value = 6010.0002
$ws.Range("L126").Value = 8208.545999999998
$ws.Range("M126").Value = -3540.0002
$ws.Range("N126").Value = -13148.546

$ws.Range("H136").Value = 201094
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 201094
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 603282
$ws.Range("N136").Value = -608382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3961.2727
$ws.Range("I7").Value = 3469.5715
$ws.Range("J7").Value = 4821.75
$ws.Range("K7").Value = 3469.5715
$ws.Range("L7").Value = 4821.75
$ws.Range("M7").Value = -3357.5715
$ws.Range("N7").Value = -5045.75

$ws.Range("H40").Value = 3389.08
$ws.Range("I40").Value = 2268.4375
$ws.Range("J40").Value = 5381.3335
$ws.Range("K40").Value = 2268.4375
$ws.Range("L40").Value = 5381.3335
$ws.Range("M40").Value = -2132.4375
$ws.Range("N40").Value = -5653.3335

$ws.Range("H50").Value = 28750.666
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 28750.666
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 28750.666
$ws.Range("N50").Value = -30024.666
$ws.Range("M50").ClearContents()

$ws.Range("H126").Value = 3961.2727
$ws.Range("I126").Value = 3469.5715
$ws.Range("J126").Value = 4821.75
$ws.Range("K126").Value = 10408.7145
$ws.Range("L126").Value = 14465.25
$ws.Range("M126").Value = -7938.7145
$ws.Range("N126").Value = -19405.25

$ws.Range("H132").Value = 3006.375
$ws.Range("I132").Value = 2615
$ws.Range("J132").Value = 4180.5
$ws.Range("K132").Value = 7845
$ws.Range("L132").Value = 12541.5
$ws.Range("M132").Value = -5315
$ws.Range("N132").Value = -17601.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2488.1765
$ws.Range("I126").Value = 1860
$ws.Range("J126").Value = 3385.5715
$ws.Range("K126").Value = 5580
$ws.Range("L126").Value = 10156.7145
$ws.Range("M126").Value = -3110
$ws.Range("N126").Value = -15539.5001
